$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '26.820.49'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '  -1.43%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.856.61'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '  -0.85%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.001'
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '  -0.04%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '304.12'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '  -1.20%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '1.000'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '  -0.09%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.5061'
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '  -0.85%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3643'
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '  -3.17%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.07160'
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '  -0.12%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.8902'
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '  -0.01%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '20.65'
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '  -0.46%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '1.858.41'
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '  -0.74%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.07437'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '  -2.08%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '92.36'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '  +3.33%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '5.224'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '  -2.25%  '
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '  -0.10%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.000008494'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '  -0.80%  '
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '  -1.06%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '1.000'
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '  -0.05%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '26.870.92'
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '  -1.43%  '
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '  -1.34%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '2.094.68'
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '  +0.34%  '
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '  -3.03%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '6.420'
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '  -1.23%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '147.19'
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '  -2.53%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '1.794'
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '  -2.75%  '
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '  -1.30%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '2.055'
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '  -3.53%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '113.04'
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = '  +0.21%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '4.630'
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = '  -2.69%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '4.658'
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = '  -1.20%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.09222'
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = '  +2.51%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.05077'
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = '  -1.63%  '
$ws.Range("B34").Value = 'ImmutableX'
$ws.Range("C34").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.7434'
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = '  -1.63%  '
$ws.Range("B35").Value = 'HuobiToken'
$ws.Range("C35").Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '2.973'
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = '  -3.95%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.143'
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = '  -2.43%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '3.246'
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = '  +7.07%  '
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '  -1.28%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.01988'
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '  -2.71%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '1.090'
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '  +1.23%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.5323'
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '  -0.74%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '118.66'
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '  +3.57%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '6.460'
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '  -3.04%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '8.371'
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '  -2.43%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.1456'
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '  -1.91%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.4641'
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '  -0.86%  '
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '  -0.18%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '9.973'
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '  -1.59%  '
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '  -0.92%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '36.75'
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '  +0.07%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '62.83'
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '  -3.64%  '
